# "Generate Report for Handoff" — refresh the localization-status report so
# that b.md's row reflects a brand-new handoff package (xliff regenerated),
# while a.md's rows are untouched.

$wb = $excel.ActiveWorkbook

$status_ready  = "Ready for handoff"
$newDate       = "2016-08-20 22:44:07"

# ---------------------------------------------------------------------
# Overview sheet: row 3 is the b.md summary row.
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("E3").Value = $status_ready
$ov.Range("F3").Value = $status_ready
$ov.Range("G3").Value = $newDate

# ---------------------------------------------------------------------
# zh-cn sheet: row 3 (b.md) gets a freshly generated handoff xliff.
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("C3").Value = $status_ready
$zh.Range("F3").Value = "'False"
$zh.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zh.Range("H3").Value = "2016-08-20 22:43:58"
$zh.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b6adfa9631904f6631db643e5edcac4a4db95af6/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/da034a09f11143b4f9aa05ec5faa63478cb3adbf/e2e/b.md."
$zh.Range("P1").ColumnWidth = 39.14

# ---------------------------------------------------------------------
# de-de sheet: row 3 (b.md) gets a freshly generated handoff xliff.
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")
$de.Range("C3").Value = $status_ready
$de.Range("F3").Value = "'False"
$de.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$de.Range("H3").Value = $newDate
$de.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b6adfa9631904f6631db643e5edcac4a4db95af6/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/da034a09f11143b4f9aa05ec5faa63478cb3adbf/e2e/b.md."
$de.Range("P1").ColumnWidth = 39.14
